# Atualização do arquivo considerando os coordenadores de compra como alçada 3
#
# "lista" sheet (NIVEL column = B) gets bumped by +1 for the rows belonging to
# the "coordenadores de compra" (ana.matsunaga, fabiano.farenzena, cyro.mello,
# abreu.marcelo, alexandre.olim) across the SUPRIMENTOS 2/3/4 groups. The row
# for abreu.marcelo in the SUPRIMENTOS 4 group (B60) additionally gets an
# underline applied to flag/highlight the change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lista")

# --- SUPRIMENTOS 2 group (rows 26-28): NIVEL 2 -> 3 ---
$ws.Range("B26").Value = 3
$ws.Range("B27").Value = 3
$ws.Range("B28").Value = 3

# --- SUPRIMENTOS 3 group (rows 41-44): NIVEL 2 -> 3, 3 -> 4 ---
$ws.Range("B41").Value = 3
$ws.Range("B42").Value = 3
$ws.Range("B43").Value = 3
$ws.Range("B44").Value = 4

# --- SUPRIMENTOS 4 group (rows 57-61): NIVEL 2 -> 3, 3 -> 4, 4 -> 5 ---
$ws.Range("B57").Value = 3
$ws.Range("B58").Value = 3
$ws.Range("B59").Value = 3
$ws.Range("B60").Value = 4
$ws.Range("B60").Font.Underline = $true
$ws.Range("B61").Value = 5

# Make "lista" the active/selected sheet & tab (it was "usuarios" before).
$ws.Activate() | Out-Null
$ws.Range("B63").Select() | Out-Null
